$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4354275
$ws.Range("J17").Value = 4768830
$ws.Range("L17").Value = 14306490
$ws.Range("N17").Value = -14306826
$ws.Range("H62").Value = 4435.5
$ws.Range("J62").Value = 6132
$ws.Range("L62").Value = 6132
$ws.Range("N62").Value = -7380
$ws.Range("H65").Value = 4435.5
$ws.Range("J65").Value = 6132
$ws.Range("L65").Value = 30660
$ws.Range("N65").Value = -36900
$ws.Range("H129").Value = 117190.24
$ws.Range("I129").Value = 252.44444
$ws.Range("J129").Value = 130858.3
$ws.Range("K129").Value = 757.33332
$ws.Range("L129").Value = 392574.9
$ws.Range("M129").Value = 4242.66668
$ws.Range("N129").Value = -402574.9
$ws.Range("H131").Value = 1740.9584
$ws.Range("I131").Value = 1268.9
$ws.Range("J131").Value = 2078.1428
$ws.Range("K131").Value = 3806.7
$ws.Range("L131").Value = 6234.428400000001
$ws.Range("M131").Value = 1233.3
$ws.Range("N131").Value = -16314.4284
$ws.Range("H132").Value = 2843.139
$ws.Range("I132").Value = 2968.0908
$ws.Range("J132").Value = 1468.6666
$ws.Range("K132").Value = 8904.2724
$ws.Range("L132").Value = 4405.9998
$ws.Range("M132").Value = -6374.2724
$ws.Range("N132").Value = -9465.9998
$ws.Range("H137").Value = 1655.7693
$ws.Range("I137").Value = 1295.7333
$ws.Range("J137").Value = 2855.889
$ws.Range("K137").Value = 3887.199900000001
$ws.Range("L137").Value = 8567.667000000001
$ws.Range("M137").Value = -1337.199900000001
$ws.Range("N137").Value = -13667.667
$ws.Range("H138").Value = 2518.1406
$ws.Range("I138").Value = 2042.1666
$ws.Range("J138").Value = 2627.9807
$ws.Range("K138").Value = 6126.4998
$ws.Range("L138").Value = 7883.9421
$ws.Range("M138").Value = -986.4997999999996
$ws.Range("N138").Value = -18163.9421

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 756.9474
$ws.Range("I2").Value = 658.8333
$ws.Range("J2").Value = 1124.875
$ws.Range("K2").Value = 658.8333
$ws.Range("L2").Value = 1124.875
$ws.Range("M2").Value = -545.8333
$ws.Range("N2").Value = -1350.875
$ws.Range("H60").Value = 13352.333
$ws.Range("I60").Value = 5000
$ws.Range("J60").Value = 17528.5
$ws.Range("K60").Value = 5000
$ws.Range("L60").Value = 17528.5
$ws.Range("M60").Value = -4267
$ws.Range("N60").Value = -18994.5
$ws.Range("H61").Value = 1610.34
$ws.Range("I61").Value = 1380.375
$ws.Range("J61").Value = 2530.2
$ws.Range("K61").Value = 1380.375
$ws.Range("L61").Value = 2530.2
$ws.Range("M61").Value = -1168.375
$ws.Range("N61").Value = -2954.2
$ws.Range("H63").Value = 1839797
$ws.Range("I63").Value = 1659.3125
$ws.Range("J63").Value = 31250000
$ws.Range("K63").Value = 1659.3125
$ws.Range("L63").Value = 31250000
$ws.Range("M63").Value = -973.3125
$ws.Range("N63").Value = -31251372
$ws.Range("H66").Value = 1839797
$ws.Range("I66").Value = 1659.3125
$ws.Range("J66").Value = 31250000
$ws.Range("K66").Value = 8296.5625
$ws.Range("L66").Value = 156250000
$ws.Range("M66").Value = -4864.5625
$ws.Range("N66").Value = -156256864
$ws.Range("H74").Value = 20001202
$ws.Range("I74").Value = 23810036
$ws.Range("K74").Value = 23810036
$ws.Range("M74").Value = -23809162
$ws.Range("H77").Value = 20001202
$ws.Range("I77").Value = 23810036
$ws.Range("K77").Value = 119050180
$ws.Range("M77").Value = -119045812
$ws.Range("H116").Value = 756.9474
$ws.Range("I116").Value = 658.8333
$ws.Range("J116").Value = 1124.875
$ws.Range("K116").Value = 658.8333
$ws.Range("L116").Value = 1124.875
$ws.Range("M116").Value = 1635.1667
$ws.Range("N116").Value = -5712.875
$ws.Range("H136").Value = 1610.34
$ws.Range("I136").Value = 1380.375
$ws.Range("J136").Value = 2530.2
$ws.Range("K136").Value = 4141.125
$ws.Range("L136").Value = 7590.599999999999
$ws.Range("M136").Value = -1591.125
$ws.Range("N136").Value = -12690.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 756.9474
$ws.Range("I3").Value = 658.8333
$ws.Range("J3").Value = 1124.875
$ws.Range("K3").Value = 658.8333
$ws.Range("L3").Value = 1124.875
$ws.Range("M3").Value = -544.8333
$ws.Range("N3").Value = -1352.875
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 37926.215
$ws.Range("I58").Value = 3262.4
$ws.Range("J58").Value = 57183.89
$ws.Range("K58").Value = 3262.4
$ws.Range("L58").Value = 57183.89
$ws.Range("M58").Value = -3059.4
$ws.Range("N58").Value = -57589.89
$ws.Range("H76").Value = 4500
$ws.Range("I76").Value = 4500
$ws.Range("K76").Value = 4500
$ws.Range("M76").Value = -4185
$ws.Range("H79").Value = 4500
$ws.Range("I79").Value = 4500
$ws.Range("K79").Value = 4500
$ws.Range("M79").Value = -3408
$ws.Range("H132").Value = 2624.8708
$ws.Range("I132").Value = 1914.0454
$ws.Range("J132").Value = 4362.4443
$ws.Range("K132").Value = 5742.1362
$ws.Range("L132").Value = 13087.3329
$ws.Range("M132").Value = -3212.1362
$ws.Range("N132").Value = -18147.3329
$ws.Range("H134").Value = 823.5714
$ws.Range("I134").Value = 764.75
$ws.Range("K134").Value = 2294.25
$ws.Range("M134").Value = 240.75
$ws.Range("H136").Value = 37926.215
$ws.Range("I136").Value = 3262.4
$ws.Range("J136").Value = 57183.89
$ws.Range("K136").Value = 9787.200000000001
$ws.Range("L136").Value = 171551.67
$ws.Range("M136").Value = -7237.200000000001
$ws.Range("N136").Value = -176651.67
$ws.Range("H141").Value = 22732.82
$ws.Range("J141").Value = 23232.236
$ws.Range("L141").Value = 23232.236
$ws.Range("N141").Value = -33592.236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1170.0222
$ws.Range("I5").Value = 929.3929000000001
$ws.Range("J5").Value = 1566.3529
$ws.Range("K5").Value = 2788.1787
$ws.Range("L5").Value = 4699.0587
$ws.Range("M5").Value = -2676.1787
$ws.Range("N5").Value = -4923.0587
$ws.Range("H52").Value = 3000
$ws.Range("J52").Value = 3000
$ws.Range("L52").Value = 9000
$ws.Range("N52").Value = -9532
$ws.Range("H92").Value = 593.1429000000001
$ws.Range("I92").Value = 350
$ws.Range("J92").Value = 633.6667
$ws.Range("K92").Value = 1050
$ws.Range("L92").Value = 1901.0001
$ws.Range("M92").Value = 198
$ws.Range("N92").Value = -4397.0001
$ws.Range("H101").Value = 9768.125
$ws.Range("J101").Value = 9768.125
$ws.Range("L101").Value = 29304.375
$ws.Range("N101").Value = -34172.375
$ws.Range("H122").Value = 919.6279
$ws.Range("J122").Value = 977.8946999999999
$ws.Range("L122").Value = 8801.052299999999
$ws.Range("N122").Value = -13701.0523
$ws.Range("H123").Value = 2898.182
$ws.Range("I123").Value = 1418
$ws.Range("J123").Value = 4131.6665
$ws.Range("K123").Value = 4254
$ws.Range("L123").Value = 12394.9995
$ws.Range("M123").Value = -1804
$ws.Range("N123").Value = -17294.9995
$ws.Range("H131").Value = 756.4299999999999
$ws.Range("I131").Value = 410
$ws.Range("J131").Value = 782.5054
$ws.Range("K131").Value = 1230
$ws.Range("L131").Value = 2347.5162
$ws.Range("M131").Value = 3810
$ws.Range("N131").Value = -12427.5162
$ws.Range("H132").Value = 726.6667
$ws.Range("I132").Value = 726.6667
$ws.Range("K132").Value = 6540.0003
$ws.Range("M132").Value = -4010.0003
$ws.Range("H135").Value = 1170.0222
$ws.Range("I135").Value = 929.3929000000001
$ws.Range("J135").Value = 1566.3529
$ws.Range("K135").Value = 8364.536100000001
$ws.Range("L135").Value = 14097.1761
$ws.Range("M135").Value = -5829.536100000001
$ws.Range("N135").Value = -19167.1761
$ws.Range("H136").Value = 3293.0833
$ws.Range("I136").Value = 906
$ws.Range("K136").Value = 2718
$ws.Range("M136").Value = 2382

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 7000
$ws.Range("J52").Value = 7000
$ws.Range("L52").Value = 7000
$ws.Range("N52").Value = -7518
$ws.Range("H53").Value = 4019.5
$ws.Range("I53").Value = 4019.5
$ws.Range("K53").Value = 4019.5
$ws.Range("M53").Value = -3388.5
$ws.Range("H102").Value = 2407.4443
$ws.Range("I102").Value = 1947.2727
$ws.Range("J102").Value = 3130.5715
$ws.Range("K102").Value = 1947.2727
$ws.Range("L102").Value = 3130.5715
$ws.Range("M102").Value = -325.2727
$ws.Range("N102").Value = -6374.5715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 547.9048
$ws.Range("I16").Value = 478
$ws.Range("K16").Value = 478
$ws.Range("M16").Value = -308
$ws.Range("H22").Value = 2824.1428
$ws.Range("I22").Value = 4093.9167
$ws.Range("K22").Value = 4093.9167
$ws.Range("M22").Value = -3798.9167
$ws.Range("H27").Value = 2824.1428
$ws.Range("I27").Value = 4093.9167
$ws.Range("K27").Value = 4093.9167
$ws.Range("M27").Value = -3986.9167
$ws.Range("H46").Value = 1740.9
$ws.Range("I46").Value = 1823.5
$ws.Range("J46").Value = 1617
$ws.Range("K46").Value = 1823.5
$ws.Range("L46").Value = 1617
$ws.Range("M46").Value = -1635.5
$ws.Range("N46").Value = -1993
$ws.Range("H122").Value = 579830.4399999999
$ws.Range("I122").Value = 936562.9399999999
$ws.Range("J122").Value = 3570.2307
$ws.Range("K122").Value = 2809688.82
$ws.Range("L122").Value = 10710.6921
$ws.Range("M122").Value = -2807238.82
$ws.Range("N122").Value = -15610.6921
$ws.Range("H130").Value = 18875
$ws.Range("J130").Value = 18875
$ws.Range("L130").Value = 18875
$ws.Range("N130").Value = -28915
$ws.Range("H136").Value = 1908.3
$ws.Range("I136").Value = 1692.9474
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 5078.8422
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -2528.8422
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2745.5557
$ws.Range("I132").Value = 1952
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 5856
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -3326
$ws.Range("N132").Value = -18057.9995
